$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mixer2-BOM")

# "Have" column (C) now matches "Qty" (B) for these rows, so the
# computed "Need" column (E, =MAX(0,B-C-D)) recalculates to 0.
$ws.Range("C6").Value2 = $ws.Range("B6").Value2
$ws.Range("C7").Value2 = $ws.Range("B7").Value2
$ws.Range("C23").Value2 = $ws.Range("B23").Value2
$ws.Range("C24").Value2 = $ws.Range("B24").Value2

# Move the active selection to A21 (single cell).
$ws.Range("A21").Select()
